# Incluida constraint de hora extra
# Set min_hours_per_week (column I) to 0 for every employee (rows 2-63) on the
# "Employees" sheet, so that overtime hours are no longer forced into the
# minimum weekly hours requirement.

$wb = $excel.ActiveWorkbook

$wsEmployees = $wb.Worksheets.Item("Employees")
$wsEmployees.Range("I2:I63").Value = 0

# Reproduce the view/selection state left on the "Employees" sheet.
$wsEmployees.Activate() | Out-Null
$wsEmployees.Range("J57").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1

# Reproduce the view/selection state on the "Parameters" sheet, which is the
# tab that ends up active/selected in the saved workbook.
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Activate() | Out-Null
$wsParameters.Range("C2").Select() | Out-Null
